# Add a "MemberId" column to the appointment sample sheet (upload route /
# sample template update), between "EmiratesId" (D) and "MobileCountryCode"
# (old E, now F), fill in the sample value, mark it as text so the numeric
# id isn't auto-formatted as a number, and leave the selection on the new
# cell the way Excel would after typing the value in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "MobileCountryCode" column (E), which
# shifts MobileCountryCode..AppointmentDateTime one column to the right.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "MemberId"

# Sample data value for the new column. Set the value first, then mark the
# cell with a text number format ("@") - same as the source workbook, which
# stores a literal numeric value but displays/treats it as text.
$ws.Range("E2").Value = 123456
$ws.Range("E2").NumberFormat = "@"

# Leave the selection on the cell that was just filled in.
[void]$ws.Range("E2").Select()
